# Update Betfair Back/Lay odds for 2025-12-31 games.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Central Coast Mariners vs Brisbane Roar
$ws.Range("I2").Value = 2.22
$ws.Range("T2").Value = 1.83
$ws.Range("AJ2").Value = 1000

# Row 3 - Al Draih vs Al Jubail
$ws.Range("N3").Value = 1.32
$ws.Range("P3").Value = 1.32
$ws.Range("T3").Value = 1.04
$ws.Range("U3").Value = 1.04

# Row 4 - Al Batin vs Abha
$ws.Range("N4").Value = 1.1
$ws.Range("T4").Value = 1.04
$ws.Range("U4").Value = 1.04

# Row 5 - Briton Ferry Llansawel vs Barry Town Utd
$ws.Range("J5").Value = 1.02
$ws.Range("L5").Value = 1.01
$ws.Range("M5").Value = 1.01
$ws.Range("N5").Value = 1.1
$ws.Range("O5").Value = 1.21
$ws.Range("P5").Value = 1.28
$ws.Range("Q5").Value = 1.21
$ws.Range("R5").Value = 1.18
$ws.Range("S5").Value = 1.21
$ws.Range("T5").Value = 1.04
$ws.Range("U5").Value = 1.04
$ws.Range("V5").Value = 1.01
$ws.Range("W5").Value = 1.01
$ws.Range("X5").Value = 990
$ws.Range("Y5").Value = 990
$ws.Range("Z5").Value = 1000
$ws.Range("AA5").Value = 1000
$ws.Range("AB5").Value = 990
$ws.Range("AC5").Value = 990
$ws.Range("AD5").Value = 990
$ws.Range("AE5").Value = 1000
$ws.Range("AF5").Value = 1000
$ws.Range("AG5").Value = 990
$ws.Range("AH5").Value = 990
$ws.Range("AI5").Value = 1000
$ws.Range("AJ5").Value = 1000
$ws.Range("AK5").Value = 1000
$ws.Range("AL5").Value = 1000
$ws.Range("AM5").Value = 1000
$ws.Range("AN5").Value = 1000
$ws.Range("AO5").Value = 1000

# Row 7 - NEOM Sports Club vs Al-Ittihad
$ws.Range("H7").Value = 1.89
$ws.Range("K7").Value = 9.8

# Row 8 - Maccabi Netanya vs Hapoel Petach Tikva
$ws.Range("G8").Value = 2.58
$ws.Range("I8").Value = 3.9
$ws.Range("J8").Value = 3.4

# Row 9 - Al-Shabab (KSA) vs Al-Quadisiya (KSA)
$ws.Range("F9").Value = 4.6
$ws.Range("G9").Value = 6.2
$ws.Range("H9").Value = 1.71
$ws.Range("I9").Value = 1.79
$ws.Range("J9").Value = 3.95
$ws.Range("K9").Value = 4.4
